# Insert a new data row at row 116 (pushes existing rows 116:216 down to 117:217)
# and populate it with a new weekly price observation, matching the rest of the
# "Hortaliza, Mercado Mayorista Lo Valledor de Santiago - Orégano" dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(116).Insert()

$ws.Cells.Item(116, 1).Value = 6
$ws.Cells.Item(116, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(116, 3).Value = "Metropolitana"
$ws.Cells.Item(116, 4).Value = 44729
$ws.Cells.Item(116, 5).Value = 13
$ws.Cells.Item(116, 6).Value = 100112029
$ws.Cells.Item(116, 7).Value = "Orégano"
$ws.Cells.Item(116, 8).Value = "Sin especificar"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 46
$ws.Cells.Item(116, 11).Value = 12000
$ws.Cells.Item(116, 12).Value = 13000
$ws.Cells.Item(116, 13).Value = 12457
$ws.Cells.Item(116, 14).Value = "$/docena de atados"
$ws.Cells.Item(116, 15).Value = "Región Metropolitana"
$ws.Cells.Item(116, 16).Value = 4152
$ws.Cells.Item(116, 17).Value = 3
$ws.Cells.Item(116, 18).Value = "Hortaliza"
